$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.254.52"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.887.20"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $ws.Range("E4").Style
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.18"
$ws.Range("D5").Style = $ws.Range("E5").Style
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $ws.Range("E6").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4678"
$ws.Range("D7").Style = $ws.Range("E7").Style
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2851"
$ws.Range("D8").Style = $ws.Range("E8").Style
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06599"
$ws.Range("D9").Style = $ws.Range("E9").Style
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.00"
$ws.Range("D10").Style = $ws.Range("E10").Style
$ws.Range("E10").Value = "  +7.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07789"
$ws.Range("D11").Style = $ws.Range("E11").Style
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.07"
$ws.Range("D12").Style = $ws.Range("E12").Style
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("D13").Value = "1.895.08"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.117"
$ws.Range("D14").Style = $ws.Range("E14").Style
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6763"
$ws.Range("D15").Style = $ws.Range("E15").Style
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "285.50"
$ws.Range("D16").Style = $ws.Range("E16").Style
$ws.Range("E16").Value = "  +12.59%  "
$ws.Range("D17").Value = "30.258.92"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D18").Style = $ws.Range("E18").Style
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.143.34"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.64"
$ws.Range("D20").Style = $ws.Range("E20").Style
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.385"
$ws.Range("D21").Style = $ws.Range("E21").Style
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007299"
$ws.Range("D22").Style = $ws.Range("E22").Style
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.189"
$ws.Range("D24").Style = $ws.Range("E24").Style
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.420"
$ws.Range("D25").Style = $ws.Range("E25").Style
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.07"
$ws.Range("D26").Style = $ws.Range("E26").Style
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.25"
$ws.Range("D27").Style = $ws.Range("E27").Style
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.996"
$ws.Range("D28").Style = $ws.Range("E28").Style
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("D29").Style = $ws.Range("E29").Style
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09739"
$ws.Range("D30").Style = $ws.Range("E30").Style
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.440"
$ws.Range("D31").Style = $ws.Range("E31").Style
$ws.Range("E31").Value = "  -5.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.164"
$ws.Range("D33").Style = $ws.Range("E33").Style
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04688"
$ws.Range("D34").Style = $ws.Range("E34").Style
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7103"
$ws.Range("D35").Style = $ws.Range("E35").Style
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.097"
$ws.Range("D36").Style = $ws.Range("E36").Style
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("D37").Style = $ws.Range("E37").Style
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01876"
$ws.Range("D38").Style = $ws.Range("E38").Style
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.686"
$ws.Range("D39").Style = $ws.Range("E39").Style
$ws.Range("E39").Value = "  +7.66%  "
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.52"
$ws.Range("D41").Style = $ws.Range("E41").Style
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.976"
$ws.Range("D42").Style = $ws.Range("E42").Style
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8690"
$ws.Range("D43").Style = $ws.Range("E43").Style
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.19"
$ws.Range("D44").Style = $ws.Range("E44").Style
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4202"
$ws.Range("D46").Style = $ws.Range("E46").Style
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "986.72"
$ws.Range("D47").Style = $ws.Range("E47").Style
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.254"
$ws.Range("D48").Style = $ws.Range("E48").Style
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.242"
$ws.Range("D49").Style = $ws.Range("E49").Style
$ws.Range("E49").Value = "  +5.79%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.08"
$ws.Range("D50").Style = $ws.Range("E50").Style
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1159"
$ws.Range("D51").Style = $ws.Range("E51").Style
$ws.Range("E51").Value = "  -3.14%  "
